$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): split the combined string into individual
# header cells, each keeping the same bold/border/centered style that
# A1 already has.
$headers = @("group1", "group2", "meandiff", "p-adj", "lower", "upper", "reject")

for ($i = 0; $i -lt $headers.Length; $i++) {
    $col = $i + 1
    $cell = $ws.Cells.Item(1, $col)
    $cell.Value2 = $headers[$i]
}

# Copy A1's formatting (bold font, thin border, centered alignment) onto
# the rest of the header row so every header cell shares style index 1.
$ws.Range("A1").Copy()
$ws.Range("B1:G1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# --- Data rows (2-7): numeric group1/group2/meandiff/p-adj/lower/upper
# plus boolean reject column.
$data = @(
    @(0, 3, 0.5953, 0.0667, -0.0285, 1.2192, $false),
    @(0, 5, 0.0288, 0.9,    -0.5951, 0.6526, $false),
    @(0, 6, -0.254, 0.6835, -0.8778, 0.3699, $false),
    @(3, 5, -0.5666, 0.0879, -1.1904, 0.0573, $false),
    @(3, 6, -0.8493, 0.0036, -1.4731, -0.2254, $true),
    @(5, 6, -0.2827, 0.616,  -0.9066, 0.3411, $false)
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $r + 2
    $rowVals = $data[$r]
    for ($c = 0; $c -lt $rowVals.Length; $c++) {
        $col = $c + 1
        $ws.Cells.Item($row, $col).Value2 = $rowVals[$c]
    }
}
